# Update the "Macro_taxonomy" sheet (CAF_RES.xlsx):
#  - Insert a new allocation row for Urban/Other (ME+MEO/LWAL) right after
#    the existing Urban/Other/MATO row, splitting that row's proportion
#    50/50 between the two taxonomies.
#  - Fill in two more previously-blank rows with the matching Rural/Other
#    split across MATO / EWV-LN / ME+MEO-LWAL.
#  - Leave the sheet on-screen at the newly edited cell, as the author did.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Macro_taxonomy")

# Insert one new row at row 16 (shifts old rows 16..62 down to 17..63,
# inheriting the formatting of the row above exactly like Excel's own
# Insert Row command).
$ws.Rows.Item(16).Insert()

# Urban / Other now splits 50/50 between MATO and ME+MEO/LWAL.
$ws.Range("D15").Value = 0.5

$ws.Range("A16").Value = "Other"
$ws.Range("B16").Value = "Urban"
$ws.Range("C16").Value = "ME+MEO/LWAL"
$ws.Range("D16").Value = 0.5

# Rural / Other (was row 24, now shifted to row 25) drops from 1 -> 0.5
# to make room for two additional taxonomies below.
$ws.Range("D25").Value = 0.5

$ws.Range("A26").Value = "Other"
$ws.Range("B26").Value = "Rural"
$ws.Range("C26").Value = "EWV/LN"
$ws.Range("D26").Value = 0.25

$ws.Range("A27").Value = "Other"
$ws.Range("B27").Value = "Rural"
$ws.Range("C27").Value = "ME+MEO/LWAL"
$ws.Range("D27").Value = 0.25

# The author left off editing with Macro_taxonomy on-screen, selection on
# the cell they'd just typed into.
$ws.Activate()
$ws.Range("D17").Select()
